$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.773.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.539.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.94"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.536.94"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.197"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.73"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.584"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.79"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000281"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.087.87"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "619.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -9.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.511.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.638.10"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.28"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.887"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.90"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.91"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.90"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.64"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.35"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.46"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.54"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "571.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.60"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.84"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.48"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.140"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0449"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.392.80"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.328"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0712"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "33.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.62"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.130"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.65"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.10%  "
